# Fixed player spawning & Docs updated
# Update the DTT Test Hour Log worksheet with the new log entry (row 7),
# corrected dates for the existing entries, and refresh the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DTT Test Hour Log")

# Correct the dates of the three already-logged entries.
$ws.Range("C4").Value = 43976
$ws.Range("C5").Value = 43977
$ws.Range("C6").Value = 43977

# Add the new log entry in row 7.
$ws.Range("A7").Value = "Implemented Shader & new Walls"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 43983
$ws.Range("D7").Value = "Implemented new walls made in blende & made Cell shaded Shader"

# Recompute the total hours formula in B30.
$excel.Calculate()

# Move the active selection to H4, matching the saved view state.
$ws.Range("H4").Select()
